$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C ("Förändrad") from 45179 to 45180 for every data row (2..496)
for ($r = 2; $r -le 496; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value = 45180
    }
}

# 2) Row 236 ("A 29118-2020") gains link formulas in columns U..Y,
#    matching the pattern used by other rows' link columns.
$ws.Range("U236").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MONSTERAS/knärot/A 29118-2020.png")'
$ws.Range("V236").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MONSTERAS/klagomål/A 29118-2020.docx")'
$ws.Range("W236").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MONSTERAS/klagomålsmail/A 29118-2020.docx")'
$ws.Range("X236").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MONSTERAS/tillsyn/A 29118-2020.docx")'
$ws.Range("Y236").Formula = '=HYPERLINK("https://klasma.github.io/Logging_MONSTERAS/tillsynsmail/A 29118-2020.docx")'
